$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.242.62'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.594.79'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.26'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.504'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0604'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.94'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0854'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '1.821.25'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '1.591.90'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.502'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.54'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '26.244.01'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.39'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.19%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0721'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.92'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.13'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.99'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.112'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0493'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').Value = '1.472.91'
$ws.Range('E32').Value = '  +3.71%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.565'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.18%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.931'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = '1.733.76'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.754'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.77'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.21%  '
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0501'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0948'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.91%  '
